# Update NATMI LR-pair data (Tnfsf13-Tnfrsf17) with newly computed TPM-derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> Resolving-Mac)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.4678846666666667
$ws.Range("H2").Value = 1.403654
$ws.Range("I2").Value = 0.03878103110832688
$ws.Range("J2").Value = 0.03913668944394218
$ws.Range("M2").Value = 0.09788266666666667
$ws.Range("N2").Value = 0.293648
$ws.Range("Q2").Value = 0.04579779886577778
$ws.Range("R2").Value = 0.412180189792
$ws.Range("S2").Value = 0.03878103110832688
$ws.Range("T2").Value = 0.03913668944394218

# Row 3 (FAPs -> Resolving-Mac)
$ws.Range("I3").Value = 0.01243343417919891
$ws.Range("J3").Value = 0.01254746040232341
$ws.Range("M3").Value = 0.09788266666666667
$ws.Range("N3").Value = 0.293648
$ws.Range("Q3").Value = 0.01468305255111111
$ws.Range("R3").Value = 0.13214747296
$ws.Range("S3").Value = 0.01243343417919891
$ws.Range("T3").Value = 0.01254746040232341

# Row 4 (Inflammatory-Mac -> Resolving-Mac)
$ws.Range("G4").Value = 5.973387666666667
$ws.Range("H4").Value = 17.920163
$ws.Range("I4").Value = 0.4951094776699161
$ws.Range("J4").Value = 0.4996500947639684
$ws.Range("M4").Value = 0.09788266666666667
$ws.Range("N4").Value = 0.293648
$ws.Range("Q4").Value = 0.5846911138471113
$ws.Range("R4").Value = 5.262220024624001
$ws.Range("S4").Value = 0.4951094776699161
$ws.Range("T4").Value = 0.4996500947639684

# Row 5 (MuSCs -> Resolving-Mac)
$ws.Range("G5").Value = 0.3289195
$ws.Range("H5").Value = 0.657839
$ws.Range("I5").Value = 0.02726278134419591
$ws.Range("J5").Value = 0.01834187103596291
$ws.Range("M5").Value = 0.09788266666666667
$ws.Range("N5").Value = 0.293648
$ws.Range("Q5").Value = 0.03219551777866667
$ws.Range("R5").Value = 0.193173106672
$ws.Range("S5").Value = 0.02726278134419591
$ws.Range("T5").Value = 0.01834187103596291

# Row 6 (Resolving-Mac -> Resolving-Mac)
$ws.Range("G6").Value = 5.144583
$ws.Range("H6").Value = 15.433749
$ws.Range("I6").Value = 0.4264132756983622
$ws.Range("J6").Value = 0.4303238843538031
$ws.Range("M6").Value = 0.09788266666666667
$ws.Range("N6").Value = 0.293648
$ws.Range("Q6").Value = 0.5035655029280001
$ws.Range("R6").Value = 4.532089526352
$ws.Range("S6").Value = 0.4264132756983622
$ws.Range("T6").Value = 0.4303238843538031
